$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Hyperlinks.Add($ws.Range("P24"), "https://www.digikey.com")
$ws.Hyperlinks.Add($ws.Range("P25"), "https://www.digikey.com")
$ws.Hyperlinks.Delete()
$ws.Range("P24:P25").Font.Underline = $false
$ws.Range("P24:P25").Font.ThemeColor = 1
